$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look like plain numbers,
# so Excel keeps them as text instead of auto-converting to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "28.116.46"
$ws.Range("E2").Value = "  +2.21%  "

$ws.Range("D3").Value = "1.909.86"
$ws.Range("E3").Value = "  +1.95%  "

$ws.Range("E4").Value = "  -1.14%  "

$ws.Range("D5").Value = "316.22"
$ws.Range("E5").Value = "  +1.05%  "

$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -1.18%  "

$ws.Range("D7").Value = "0.4836"
$ws.Range("E7").Value = "  +1.00%  "

$ws.Range("D8").Value = "0.3820"
$ws.Range("E8").Value = "  +1.04%  "

$ws.Range("D9").Value = "0.07355"
$ws.Range("E9").Value = "  -0.32%  "

$ws.Range("D10").Value = "0.9344"
$ws.Range("E10").Value = "  -0.41%  "

$ws.Range("D11").Value = "20.79"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "0.07789"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").Value = "1.942.36"
$ws.Range("E13").Value = "  +3.69%  "

$ws.Range("D14").Value = "5.511"
$ws.Range("E14").Value = "  +1.14%  "

$ws.Range("D15").Value = "6.638"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").Value = "91.51"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  -1.10%  "

$ws.Range("D18").Value = "0.000008833"
$ws.Range("E18").Value = "  -0.89%  "

$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").Value = "28.156.06"
$ws.Range("E20").Value = "  +2.24%  "

$ws.Range("D21").Value = "14.86"
$ws.Range("E21").Value = "  -0.43%  "

$ws.Range("D22").Value = "5.160"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("D23").Value = "2.144.93"
$ws.Range("E23").Value = "  +1.47%  "

$ws.Range("D24").Value = "10.90"
$ws.Range("E24").Value = "  +1.49%  "

$ws.Range("D25").Value = "156.42"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("D26").Value = "1.913"
$ws.Range("E26").Value = "  -2.61%  "

$ws.Range("D27").Value = "18.55"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").Value = "2.117"
$ws.Range("E28").Value = "  +4.78%  "

$ws.Range("D29").Value = "116.30"
$ws.Range("E29").Value = "  +0.22%  "

$ws.Range("D30").Value = "4.950"
$ws.Range("E30").Value = "  -0.95%  "

$ws.Range("D31").Value = "0.08925"
$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("D32").Value = "3.340"
$ws.Range("E32").Value = "  -0.34%  "

$ws.Range("D33").Value = "1.253"
$ws.Range("E33").Value = "  +2.74%  "

$ws.Range("D34").Value = "0.7692"
$ws.Range("E34").Value = "  +2.03%  "

$ws.Range("D35").Value = "4.679"
$ws.Range("E35").Value = "  +1.40%  "

$ws.Range("D36").Value = "2.614"
$ws.Range("E36").Value = "  -3.25%  "

$ws.Range("D37").Value = "0.02053"
$ws.Range("E37").Value = "  -0.57%  "

$ws.Range("D38").Value = "1.102"
$ws.Range("E38").Value = "  -1.47%  "

$ws.Range("D39").Value = "0.05311"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").Value = "0.5497"
$ws.Range("E40").Value = "  +2.47%  "

$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").Value = "7.013"
$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("D44").Value = "8.465"
$ws.Range("E44").Value = "  +0.35%  "

$ws.Range("D45").Value = "10.74"
$ws.Range("E45").Value = "  +1.30%  "

$ws.Range("D46").Value = "0.4834"
$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").Value = "107.34"
$ws.Range("E47").Value = "  +4.21%  "

$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -1.30%  "

$ws.Range("D49").Value = "1.658"
$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("D50").Value = "68.37"
$ws.Range("E50").Value = "  +1.46%  "

$ws.Range("E51").Value = "  +0.26%  "

